$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; this shifts the existing rows 18-97
# down to 19-98 (their contents/styles move with them).
$ws.Rows.Item(18).Insert()

# Fill the newly inserted row 18 with the new "book manager" account
# (matches the commit message: "them tai khoan quan ly sach" = "add book
# manager account").
$ws.Range("A18").Value = "Quản lý sách"
$ws.Range("B18").Value = "sba_manager"
$ws.Range("C18").Value = 22
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 2

# Move the selection to match the post-edit cursor position recorded in
# the saved workbook.
$ws.Range("F18").Select()
